$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no movie was selected for Friday. If further action is required or if you have more questions, feel free to ask!`n"
$ws.Range("C3").Value = "MSG: None`n`nMSG: The decision process concluded without selecting a movie. No action will be taken.`n"
$ws.Range("C4").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision.`"`n"
$ws.Range("C5").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no choice of a movie is possible without further discussion.`n"
$ws.Range("C6").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie was selected for Friday.`n"
$ws.Range("C7").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Oppenheimer`" to be shown on Friday.`n"
$ws.Range("C8").Value = "MSG: None`n`nMSG: The committee did not reach a decision about which movie to show on Friday.`n"
$ws.Range("C9").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("C10").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Oppenheimer`" has been recorded.`n"
$ws.Range("C11").Value = "MSG: None`n`nMSG: The function call for ``no_decision()`` has been executed, indicating that the committee did not agree on a movie for Friday.`n"
$ws.Range("C12").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday was not reached.`n"
$ws.Range("C13").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision was made regarding the movie to be shown on Friday.`n"
$ws.Range("C14").Value = "MSG: None`n`nMSG: The decision has been recorded, and the rights for `"Barbie`" will be acquired for the show on Friday.`n"
$ws.Range("C15").Value = "MSG: None`n`nMSG: The decision has been recorded. No movie was selected for Friday.`n"
$ws.Range("C16").Value = "MSG: None`n`nMSG: The decision has been made that no movie will be selected for Friday.`n"
$ws.Range("C17").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie was chosen for screening on Friday.`n"
$ws.Range("C18").Value = "MSG: None`n`nMSG: The rights for both movies have been acquired successfully.`n"
$ws.Range("C19").Value = "MSG: None`n`nMSG: The decision has been made to show `"Barbie`" on Friday.`n"
$ws.Range("C20").Value = "MSG: None`n`nMSG: The function has been successfully called to indicate that no decision was made regarding the movie selection.`n"
$ws.Range("C21").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights to both movies.`n"
$ws.Range("D21").Value = "both_movies, "
$ws.Range("C22").Value = "MSG: None`n`nMSG: The decision-making process did not lead to a consensus on the movie to be shown. Thus, I have recorded that no decision was made regarding Friday's movie.`n"
$ws.Range("C23").Value = "MSG: None`n`nMSG: The decision has been recorded, and `"Oppenheimer`" will be acquired for the movie being shown on Friday.`n"
$ws.Range("C24").Value = "MSG: None`n`nMSG: The decision has been recorded as that no movie was selected in this meeting.`n"
$ws.Range("C25").Value = "MSG: None`n`nMSG: The rights to both movies have been successfully acquired for showing on Friday.`n"
$ws.Range("C26").Value = "MSG: None`n`nMSG: The decision-making committee has not reached an agreement on which movie to select.`n"
$ws.Range("C27").Value = "MSG: None`n`nMSG: The decision has been recorded, and the rights for `"Barbie`" will be acquired for the movie shown on Friday.`n"
$ws.Range("C28").Value = "MSG: None`n`nMSG: The decision-making process did not result in a choice for the movie to be shown on Friday.`n"
$ws.Range("C29").Value = "MSG: None`n`nMSG: No decision was made regarding the movie selection for Friday.`n"
$ws.Range("C31").Value = "MSG: None`n`nMSG: I have recorded that no decision was made regarding the movie selection for Friday, as there was no consensus reached in the conversation.`n"
$ws.Range("C32").Value = "MSG: None`n`nMSG: The decision to acquire the rights for both movies has been successfully recorded.`n"
$ws.Range("C33").Value = "MSG: None`n`nMSG: The decision about which movie to play on Friday cannot be made.`n"
$ws.Range("C34").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday is unclear, so no action will be taken.`n"
$ws.Range("D34").Value = "no_decision, "
$ws.Range("C35").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday was left unresolved.`n"
$ws.Range("C36").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday was not reached.`n"
$ws.Range("C37").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for both movies.`n"
$ws.Range("C38").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision.`"`n"
$ws.Range("C39").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision.`"`n"
$ws.Range("C40").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been recorded successfully.`n"
$ws.Range("C41").Value = "MSG: None`n`nMSG: The decision regarding Friday's movie has been recorded as `"no decision.`"`n"
$ws.Range("D41").Value = "no_decision, "
$ws.Range("C42").Value = "MSG: None`n`nMSG: The decision-making process concluded without a definitive plan for Friday's movie, resulting in no decision being made.`n"
$ws.Range("C43").Value = "MSG: None`n`nMSG: No decision was made regarding the movie selection.`n"
$ws.Range("C44").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Oppenheimer`" has been recorded successfully.`n"
$ws.Range("C45").Value = "MSG: None`n`nMSG: The decision regarding the movie to show on Friday has resulted in no final choice being made.`n"
$ws.Range("D45").Value = "no_decision, "
$ws.Range("C46").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday was not reached.`n"
$ws.Range("D46").Value = "no_decision, "
$ws.Range("C47").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights to `"Barbie.`"`n"
$ws.Range("C48").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie.`"`n"
$ws.Range("C49").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie has been selected for Friday.`n"
$ws.Range("C50").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday remains unresolved, so no acquisition will be made.`n"
$ws.Range("C51").Value = "MSG: None`n`nMSG: The decision process did not result in a selection for Friday's movie.`n"
$ws.Range("C52").Value = "MSG: None`n`nMSG: The decision has been recorded as a no decision regarding the movie to be shown on Friday.`n"
$ws.Range("C53").Value = "MSG: None`n`nMSG: The committee did not arrive at a decision regarding which movie to show on Friday, so the outcome is a no decision.`n"
$ws.Range("C54").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`" regarding the movie to be shown on Friday.`n"
$ws.Range("C55").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie was selected in the meeting.`n"
$ws.Range("C56").Value = "MSG: None`n`nMSG: The decision has been recorded as no agreement was reached on which movie to show on Friday.`n"
$ws.Range("C57").Value = "MSG: None`n`nMSG: The decision regarding the movie to be shown on Friday resulted in no selection being made.`n"
$ws.Range("D57").Value = "no_decision, "
$ws.Range("C58").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no_decision`" since the committee did not reach an agreement on the movie to be shown on Friday.`n"
